# Auto-generated edit script: re-applies scraped market-price recalculation
# to the Tonberry_Profits workbook (per scheduled-runner commit).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 24.333334
$ws.Range("I11").Value = 24.333334
$ws.Range("K11").Value = 24.333334
$ws.Range("M11").Value = 115.666666
$ws.Range("H18").Value = 12460.654
$ws.Range("I18").Value = 7356.9287
$ws.Range("J18").Value = 18415
$ws.Range("K18").Value = 7356.9287
$ws.Range("L18").Value = 18415
$ws.Range("M18").Value = -7072.9287
$ws.Range("N18").Value = -18983
$ws.Range("H19").Value = 1428.2963
$ws.Range("I19").Value = 1234.2307
$ws.Range("K19").Value = 1234.2307
$ws.Range("M19").Value = -1059.2307
$ws.Range("H32").Value = 3959.1667
$ws.Range("J32").Value = 3959.1667
$ws.Range("L32").Value = 3959.1667
$ws.Range("N32").Value = -4611.1667
$ws.Range("H43").Value = 1375.75
$ws.Range("I43").Value = 856
$ws.Range("J43").Value = 1549
$ws.Range("K43").Value = 856
$ws.Range("L43").Value = 1549
$ws.Range("M43").Value = -787
$ws.Range("N43").Value = -1687
$ws.Range("H51").Value = 4190.5454
$ws.Range("I51").Value = 2998.5
$ws.Range("J51").Value = 4455.4443
$ws.Range("K51").Value = 2998.5
$ws.Range("L51").Value = 4455.4443
$ws.Range("M51").Value = -2514.5
$ws.Range("N51").Value = -5423.4443
$ws.Range("H82").Value = 2740.1667
$ws.Range("I82").Value = 2740.1667
$ws.Range("K82").Value = 8220.500100000001
$ws.Range("M82").Value = -7814.500100000001
$ws.Range("H85").Value = 2740.1667
$ws.Range("I85").Value = 2740.1667
$ws.Range("K85").Value = 8220.500100000001
$ws.Range("M85").Value = -6816.500100000001
$ws.Range("H112").Value = 2510.4443
$ws.Range("J112").Value = 2510.4443
$ws.Range("L112").Value = 7531.3329
$ws.Range("N112").Value = -9747.332900000001
$ws.Range("H116").Value = 9422.23
$ws.Range("I116").Value = 15869.857
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 15869.857
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = -12427.857
$ws.Range("N116").Value = -8784
$ws.Range("H129").Value = 885.4648
$ws.Range("J129").Value = 879.6984
$ws.Range("L129").Value = 2639.0952
$ws.Range("N129").Value = -12639.0952
$ws.Range("H132").Value = 935.5714
$ws.Range("I132").Value = 875.2083
$ws.Range("J132").Value = 1297.75
$ws.Range("K132").Value = 2625.6249
$ws.Range("L132").Value = 3893.25
$ws.Range("M132").Value = -95.6248999999998
$ws.Range("N132").Value = -8953.25
$ws.Range("H137").Value = 2108.6538
$ws.Range("I137").Value = 1163
$ws.Range("K137").Value = 3489
$ws.Range("M137").Value = -939
$ws.Range("H138").Value = 2306.2927
$ws.Range("J138").Value = 2259.8572
$ws.Range("L138").Value = 6779.571599999999
$ws.Range("N138").Value = -17059.5716

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2464.5
$ws.Range("I32").Value = 1750.5646
$ws.Range("K32").Value = 1750.5646
$ws.Range("M32").Value = -1463.5646
$ws.Range("H35").Value = 3668
$ws.Range("I35").Value = 3668
$ws.Range("K35").Value = 3668
$ws.Range("M35").Value = -3262
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H61").Value = 3412.818
$ws.Range("I61").Value = 2305.4666
$ws.Range("K61").Value = 2305.4666
$ws.Range("M61").Value = -2093.4666
$ws.Range("H110").Value = 816.4545000000001
$ws.Range("I110").Value = 808.1
$ws.Range("K110").Value = 808.1
$ws.Range("M110").Value = 1236.9
$ws.Range("H132").Value = 2096.9302
$ws.Range("I132").Value = 1837.6666
$ws.Range("K132").Value = 5512.9998
$ws.Range("M132").Value = -2982.9998
$ws.Range("H136").Value = 3412.818
$ws.Range("I136").Value = 2305.4666
$ws.Range("K136").Value = 6916.399800000001
$ws.Range("M136").Value = -4366.399800000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1824.8636
$ws.Range("J20").Value = 1762.125
$ws.Range("L20").Value = 1762.125
$ws.Range("N20").Value = -2256.125
$ws.Range("H36").Value = 7312.3335
$ws.Range("I36").Value = 968.5
$ws.Range("K36").Value = 968.5
$ws.Range("M36").Value = -434.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1321.6052
$ws.Range("I31").Value = 829.1875
$ws.Range("J31").Value = 1679.7273
$ws.Range("K31").Value = 829.1875
$ws.Range("L31").Value = 1679.7273
$ws.Range("M31").Value = -534.1875
$ws.Range("N31").Value = -2269.7273
$ws.Range("H34").Value = 1321.6052
$ws.Range("I34").Value = 829.1875
$ws.Range("J34").Value = 1679.7273
$ws.Range("K34").Value = 829.1875
$ws.Range("L34").Value = 1679.7273
$ws.Range("M34").Value = -627.1875
$ws.Range("N34").Value = -2083.7273
$ws.Range("H58").Value = 3346115.5
$ws.Range("J58").Value = 1389.8334
$ws.Range("L58").Value = 1389.8334
$ws.Range("N58").Value = -1795.8334
$ws.Range("H134").Value = 1915.05
$ws.Range("I134").Value = 1465.75
$ws.Range("K134").Value = 4397.25
$ws.Range("M134").Value = -1862.25
$ws.Range("H136").Value = 3346115.5
$ws.Range("J136").Value = 1389.8334
$ws.Range("L136").Value = 4169.5002
$ws.Range("N136").Value = -9269.5002
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 498.8
$ws.Range("J34").Value = 1998
$ws.Range("L34").Value = 5994
$ws.Range("N34").Value = -6162
$ws.Range("H39").Value = 1849
$ws.Range("J39").Value = 2999
$ws.Range("L39").Value = 8997
$ws.Range("N39").Value = -9585
$ws.Range("H55").Value = 3498.75
$ws.Range("J55").Value = 4997.5
$ws.Range("L55").Value = 14992.5
$ws.Range("N55").Value = -15346.5
$ws.Range("H62").Value = 3924.375
$ws.Range("J62").Value = 3924.375
$ws.Range("L62").Value = 11773.125
$ws.Range("N62").Value = -13145.125
$ws.Range("H65").Value = 3924.375
$ws.Range("J65").Value = 3924.375
$ws.Range("L65").Value = 35319.375
$ws.Range("N65").Value = -42183.375
$ws.Range("H104").Value = 5370.8
$ws.Range("J104").Value = 5562
$ws.Range("L104").Value = 16686
$ws.Range("N104").Value = -21928
$ws.Range("H129").Value = 48859.8
$ws.Range("I129").Value = 676
$ws.Range("J129").Value = 91020.625
$ws.Range("K129").Value = 2028
$ws.Range("L129").Value = 273061.875
$ws.Range("M129").Value = 2972
$ws.Range("N129").Value = -283061.875
$ws.Range("H131").Value = 8211002.5
$ws.Range("I131").Value = 125000420
$ws.Range("J131").Value = 15253.895
$ws.Range("K131").Value = 375001260
$ws.Range("L131").Value = 45761.685
$ws.Range("M131").Value = -374996220
$ws.Range("N131").Value = -55841.685
$ws.Range("H140").Value = 3354.7727
$ws.Range("J140").Value = 4602.467
$ws.Range("L140").Value = 13807.401
$ws.Range("N140").Value = -24167.401

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9375.643
$ws.Range("J80").Value = 10634.167
$ws.Range("L80").Value = 10634.167
$ws.Range("N80").Value = -12630.167
$ws.Range("H83").Value = 9375.643
$ws.Range("J83").Value = 10634.167
$ws.Range("L83").Value = 53170.835
$ws.Range("N83").Value = -63154.835
$ws.Range("H132").Value = 1014097.2
$ws.Range("I132").Value = 1480690.6
$ws.Range("J132").Value = 3144.75
$ws.Range("K132").Value = 4442071.800000001
$ws.Range("L132").Value = 9434.25
$ws.Range("M132").Value = -4439541.800000001
$ws.Range("N132").Value = -14494.25

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H132").Value = 4253.4546
$ws.Range("I132").Value = 1390
$ws.Range("K132").Value = 4170
$ws.Range("M132").Value = -1640

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H96").Value = 4461.5386
$ws.Range("J96").Value = 4461.5386
$ws.Range("L96").Value = 4461.5386
$ws.Range("N96").Value = -7207.5386
$ws.Range("H126").Value = 4259.0435
$ws.Range("I126").Value = 3317.111
$ws.Range("K126").Value = 9951.332999999999
$ws.Range("M126").Value = -7481.332999999999
$ws.Range("H128").Value = 29999.5
$ws.Range("J128").Value = 29999.5
$ws.Range("L128").Value = 29999.5
$ws.Range("N128").Value = -39959.5
$ws.Range("H132").Value = 1873.6061
$ws.Range("I132").Value = 1474.1818
$ws.Range("J132").Value = 2672.4546
$ws.Range("K132").Value = 4422.5454
$ws.Range("L132").Value = 8017.3638
$ws.Range("M132").Value = -1892.5454
$ws.Range("N132").Value = -13077.3638

